# Bot example for RPA.Excel.Application
# Writes a single requirement line into the blank dev-data workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "pywin32==302"
